$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 17019
$ws1.Range("F17").Value = 120
$ws1.Range("F24").Value = 7003
$ws1.Range("F28").Value = 21
$ws1.Range("F35").Value = 4962

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 17019
$ws4.Range("F17").Value = 120
$ws4.Range("F25").Value = 7003
$ws4.Range("F29").Value = 21
$ws4.Range("F37").Value = 4962
